# Update computed line-flow results (pl_mw) for the 380 kV case
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 3.768577808579209
$ws.Range("C2").Value = 0.2696893970263829
$ws.Range("D2").Value = 0.01177338621349833
$ws.Range("F2").Value = 4.561091461841045
$ws.Range("G2").Value = 0.002645219141746815
$ws.Range("J2").Value = 0.1334494828089401
$ws.Range("L2").Value = 0.3298031688368681
$ws.Range("M2").Value = 0.702488960043226
$ws.Range("N2").Value = 3.148878555820716

$ws.Range("B3").Value = 3.671253686709235
$ws.Range("C3").Value = 0.2449592784458332
$ws.Range("D3").Value = 0.01038945082990494
$ws.Range("F3").Value = 4.53715470808099
$ws.Range("G3").Value = 0.002651376000822173
$ws.Range("J3").Value = 0.1334652963854541
$ws.Range("L3").Value = 0.3288900371184269
$ws.Range("M3").Value = 0.6895179128371325
$ws.Range("N3").Value = 3.164640919489656

$ws.Range("B4").Value = 3.613753226221718
$ws.Range("C4").Value = 0.2299275255790008
$ws.Range("D4").Value = 0.009537214601323285
$ws.Range("F4").Value = 4.524608119840096
$ws.Range("G4").Value = 0.0026553551481452
$ws.Range("J4").Value = 0.1334769510082596
$ws.Range("L4").Value = 0.3284653301455549
$ws.Range("M4").Value = 0.6819487420823123
$ws.Range("N4").Value = 3.175130538977314

$ws.Range("B5").Value = 3.590887279277922
$ws.Range("C5").Value = 0.2238397629052997
$ws.Range("D5").Value = 0.009189208388924897
$ws.Range("F5").Value = 4.520034640097151
$ws.Range("G5").Value = 0.002657026851696079
$ws.Range("J5").Value = 0.1334821882214179
$ws.Range("L5").Value = 0.3283264787878224
$ws.Range("M5").Value = 0.6789634268715332
$ws.Range("N5").Value = 3.179608710100695

$ws.Range("B6").Value = 3.587124543795596
$ws.Range("C6").Value = 0.2228311605999806
$ws.Range("D6").Value = 0.009131376659702539
$ws.Range("F6").Value = 4.519307751303714
$ws.Range("G6").Value = 0.002657307471940135
$ws.Range("J6").Value = 0.1334830872871056
$ws.Range("L6").Value = 0.3283054904166605
$ws.Range("M6").Value = 0.678473703204844
$ws.Range("N6").Value = 3.180364587852324

$ws.Range("B7").Value = 3.613442558598308
$ws.Range("C7").Value = 0.2298452715986627
$ws.Range("D7").Value = 0.009532524266894882
$ws.Range("F7").Value = 4.524544258405214
$ws.Range("G7").Value = 0.002655377489954605
$ws.Range("J7").Value = 0.1334770196654347
$ws.Range("L7").Value = 0.3284633189476125
$ws.Range("M7").Value = 0.6819080797107588
$ws.Range("N7").Value = 3.175190109602511

$ws.Range("B8").Value = 3.734551024899872
$ws.Range("C8").Value = 0.2611304104252667
$ws.Range("D8").Value = 0.01129667085753283
$ws.Range("F8").Value = 4.552390799994399
$ws.Range("G8").Value = 0.00264730087003362
$ws.Range("J8").Value = 0.1334545308801696
$ws.Range("L8").Value = 0.3294601203402152
$ws.Range("M8").Value = 0.6979344293532748
$ws.Range("N8").Value = 3.154144707395716

$ws.Range("B9").Value = 3.990042169599633
$ws.Range("C9").Value = 0.3237214297756452
$ws.Range("D9").Value = 0.01474037752573309
$ws.Range("F9").Value = 4.624136557265729
$ws.Range("G9").Value = 0.00263303207267862
$ws.Range("J9").Value = 0.1334259235310187
$ws.Range("L9").Value = 0.3324927434640372
$ws.Range("M9").Value = 0.7325075247457278
$ws.Range("N9").Value = 3.119335241197462

$ws.Range("B10").Value = 4.188875820627004
$ws.Range("C10").Value = 0.370509917794152
$ws.Range("D10").Value = 0.01726660742614428
$ws.Range("F10").Value = 4.687411855185246
$ws.Range("G10").Value = 0.002623494280428629
$ws.Range("J10").Value = 0.1334144409488562
$ws.Range("L10").Value = 0.3353775311626492
$ws.Range("M10").Value = 0.7598446350984247
$ws.Range("N10").Value = 3.097729699562933

$ws.Range("B11").Value = 4.281781327981037
$ws.Range("C11").Value = 0.3919799786056046
$ws.Range("D11").Value = 0.01841628057594846
$ws.Range("F11").Value = 4.718517454515364
$ws.Range("G11").Value = 0.002619358186925448
$ws.Range("J11").Value = 0.1334113085912367
$ws.Range("L11").Value = 0.3368324955508655
$ws.Range("M11").Value = 0.7727058530115585
$ws.Range("N11").Value = 3.088769350219678

$ws.Range("B12").Value = 4.317317695224233
$ws.Range("C12").Value = 0.4001376006352189
$ws.Range("D12").Value = 0.01885180054890867
$ws.Range("F12").Value = 4.73063210979501
$ws.Range("G12").Value = 0.002617820918564568
$ws.Range("J12").Value = 0.1334104248980799
$ws.Range("L12").Value = 0.3374039565604079
$ws.Range("M12").Value = 0.777637532339142
$ws.Range("N12").Value = 3.08550171707671

$ws.Range("B13").Value = 4.309648477975315
$ws.Range("C13").Value = 0.3983794824859501
$ws.Range("D13").Value = 0.01875799472696116
$ws.Range("F13").Value = 4.728008044729506
$ws.Range("G13").Value = 0.002618150710437609
$ws.Range("J13").Value = 0.1334106017362355
$ws.Range("L13").Value = 0.3372799708185283
$ws.Range("M13").Value = 0.7765726715916372
$ws.Range("N13").Value = 3.086199870589894

$ws.Range("B14").Value = 4.284697798363823
$ws.Range("C14").Value = 0.3926505585649238
$ws.Range("D14").Value = 0.01845210723686819
$ws.Range("F14").Value = 4.719507397137448
$ws.Range("G14").Value = 0.002619231134969896
$ws.Range("J14").Value = 0.1334112298169536
$ws.Range("L14").Value = 0.3368790992791304
$ws.Range("M14").Value = 0.7731103529831742
$ws.Range("N14").Value = 3.088498002136816

$ws.Range("B15").Value = 4.269461088293951
$ws.Range("C15").Value = 0.3891450123808795
$ws.Range("D15").Value = 0.01826476647187292
$ws.Range("F15").Value = 4.714344270317099
$ws.Range("G15").Value = 0.002619896695248939
$ws.Range("J15").Value = 0.1334116539775136
$ws.Range("L15").Value = 0.3366362228117197
$ws.Range("M15").Value = 0.7709975881910864
$ws.Range("N15").Value = 3.089922030868422

$ws.Range("B16").Value = 4.182853801891213
$ws.Range("C16").Value = 0.3691105931476955
$ws.Range("D16").Value = 0.01719149129905162
$ws.Range("F16").Value = 4.685425898362524
$ws.Range("G16").Value = 0.002623768648167992
$ws.Range("J16").Value = 0.1334146878837785
$ws.Range("L16").Value = 0.3352853150246773
$ws.Range("M16").Value = 0.7590127052061177
$ws.Range("N16").Value = 3.098332801193706

$ws.Range("B17").Value = 4.130353267325916
$ws.Range("C17").Value = 0.3568681695470559
$ws.Range("D17").Value = 0.01653326148274914
$ws.Range("F17").Value = 4.668281258598086
$ws.Range("G17").Value = 0.002626195762328438
$ws.Range("J17").Value = 0.1334170859067099
$ws.Range("L17").Value = 0.3344931029606784
$ws.Range("M17").Value = 0.7517694974639539
$ws.Range("N17").Value = 3.10371532196848

$ws.Range("B18").Value = 4.100387307278709
$ws.Range("C18").Value = 0.3498440937449345
$ws.Range("D18").Value = 0.01615470440039957
$ws.Range("F18").Value = 4.658638499700089
$ws.Range("G18").Value = 0.002627610862958801
$ws.Range("J18").Value = 0.1334186619029456
$ws.Range("L18").Value = 0.3340508722532007
$ws.Range("M18").Value = 0.7476434245131287
$ws.Range("N18").Value = 3.106892879301142

$ws.Range("B19").Value = 4.090280954620937
$ws.Range("C19").Value = 0.3474688380160273
$ws.Range("D19").Value = 0.01602653535345411
$ws.Range("F19").Value = 4.655411079720807
$ws.Range("G19").Value = 0.002628093275271534
$ws.Range("J19").Value = 0.1334192292420591
$ws.Range("L19").Value = 0.3339034474712648
$ws.Range("M19").Value = 0.74625327491718
$ws.Range("N19").Value = 3.107982754355106

$ws.Range("B20").Value = 4.135918120988549
$ws.Range("C20").Value = 0.3581695864211838
$ws.Range("D20").Value = 0.01660332641753115
$ws.Range("F20").Value = 4.670083719665456
$ws.Range("G20").Value = 0.002625935417664094
$ws.Range("J20").Value = 0.1334168102609254
$ws.Range("L20").Value = 0.3345760455932876
$ws.Range("M20").Value = 0.7525364047877403
$ws.Range("N20").Value = 3.103133886051793

$ws.Range("B21").Value = 4.292016764849905
$ws.Range("C21").Value = 0.3943325340475781
$ws.Range("D21").Value = 0.01854194860375458
$ws.Range("F21").Value = 4.721995119419034
$ws.Range("G21").Value = 0.002618913002795193
$ws.Range("J21").Value = 0.1334110371111628
$ws.Range("L21").Value = 0.3369962887934719
$ws.Range("M21").Value = 0.7741256511433576
$ws.Range("N21").Value = 3.087819575322726

$ws.Range("B22").Value = 4.396106846416444
$ws.Range("C22").Value = 0.4181269965581009
$ws.Range("D22").Value = 0.01980994927880175
$ws.Range("F22").Value = 4.75787946498221
$ws.Range("G22").Value = 0.002614492293399795
$ws.Range("J22").Value = 0.1334090278760742
$ws.Range("L22").Value = 0.3386975267845855
$ws.Range("M22").Value = 0.7885935642188144
$ws.Range("N22").Value = 3.07854234180931

$ws.Range("B23").Value = 4.340361864455758
$ws.Range("C23").Value = 0.4054125913185658
$ws.Range("D23").Value = 0.01913307040790357
$ws.Range("F23").Value = 4.738547610466753
$ws.Range("G23").Value = 0.00261683631496505
$ws.Range("J23").Value = 0.133409938250554
$ws.Range("L23").Value = 0.3377786176456254
$ws.Range("M23").Value = 0.7808389198204537
$ws.Range("N23").Value = 3.083426630367725

$ws.Range("B24").Value = 4.133401573969934
$ws.Range("C24").Value = 0.3575811714967472
$ws.Range("D24").Value = 0.01657165047883069
$ws.Range("F24").Value = 4.669268160832814
$ws.Range("G24").Value = 0.002626053058040201
$ws.Range("J24").Value = 0.1334169342655445
$ws.Range("L24").Value = 0.3345385060429322
$ws.Range("M24").Value = 0.7521895671955576
$ws.Range("N24").Value = 3.103396494518321

$ws.Range("B25").Value = 3.918982718575478
$ws.Range("C25").Value = 0.3066512961679848
$ws.Range("D25").Value = 0.01380986126574157
$ws.Range("F25").Value = 4.602880806601831
$ws.Range("G25").Value = 0.002636725310052454
$ws.Range("J25").Value = 0.1334319948340266
$ws.Range("L25").Value = 0.3315569899044348
$ws.Range("M25").Value = 0.7228159417075517
$ws.Range("N25").Value = 3.128057368837659
